$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Core content edit -------------------------------------------------
# Row 4 held the per-user detail-row template: A4=${user.userCode},
# B4=${user.userName}, C4=${user.userAge}, D4=${user.userDept},
# E4=${user.userSalary}. The "${user.userName}" placeholder in B4 was
# removed (cell cleared), leaving the other cells untouched.
$ws.Range("B4").ClearContents()

# --- Cosmetic comment-formatting edit -----------------------------------
# The bold attribute was removed from the rich-text run used for the
# cell comments, leaving the rest of the run formatting (9pt Tahoma,
# indexed color 81) intact.
$cmt = $ws.Range("A4").Comment
if ($cmt -ne $null) {
    $len = $cmt.Text().Length
    $cmt.Shape.TextFrame.Characters(1, $len).Font.Bold = $false
}

# --- Selection left behind by the editing user --------------------------
$ws.Range("B2").Select() | Out-Null
